$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "Förändrad" (changed) date in column C for all existing data rows (2-499)
#    from 2023-09-21 (45190) to 2023-09-23 (45192).
$ws.Range("C2:C499").Value = 45192

# 2. Row 499 previously had no explicit row height; it now gets one (15pt, custom height),
#    matching the style used by the newly appended rows.
$ws.Rows.Item(499).RowHeight = 15

# 3. Append new row 500: "A 44566-2023"
$ws.Rows.Item(500).RowHeight = 15
$ws.Cells.Item(500, 1).Value = "A 44566-2023"
$ws.Cells.Item(500, 2).Value = 45189
$ws.Cells.Item(500, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(500, 3).Value = 45192
$ws.Cells.Item(500, 3).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(500, 4).Value = "JÖNKÖPINGS LÄN"
$ws.Cells.Item(500, 5).Value = "NÄSSJÖ"
$ws.Cells.Item(500, 7).Value = 11.8
$ws.Cells.Item(500, 8).Value = 0
$ws.Cells.Item(500, 9).Value = 0
$ws.Cells.Item(500, 10).Value = 0
$ws.Cells.Item(500, 11).Value = 0
$ws.Cells.Item(500, 12).Value = 0
$ws.Cells.Item(500, 13).Value = 0
$ws.Cells.Item(500, 14).Value = 0
$ws.Cells.Item(500, 15).Value = 0
$ws.Cells.Item(500, 16).Value = 0
$ws.Cells.Item(500, 17).Value = 0
$ws.Cells.Item(500, 18).WrapText = $true

# 4. Append new row 501: "A 44766-2023" (no explicit row height set for this row)
$ws.Cells.Item(501, 1).Value = "A 44766-2023"
$ws.Cells.Item(501, 2).Value = 45190
$ws.Cells.Item(501, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(501, 3).Value = 45192
$ws.Cells.Item(501, 3).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item(501, 4).Value = "JÖNKÖPINGS LÄN"
$ws.Cells.Item(501, 5).Value = "NÄSSJÖ"
$ws.Cells.Item(501, 7).Value = 2.3
$ws.Cells.Item(501, 8).Value = 0
$ws.Cells.Item(501, 9).Value = 0
$ws.Cells.Item(501, 10).Value = 0
$ws.Cells.Item(501, 11).Value = 0
$ws.Cells.Item(501, 12).Value = 0
$ws.Cells.Item(501, 13).Value = 0
$ws.Cells.Item(501, 14).Value = 0
$ws.Cells.Item(501, 15).Value = 0
$ws.Cells.Item(501, 16).Value = 0
$ws.Cells.Item(501, 17).Value = 0
$ws.Cells.Item(501, 18).WrapText = $true
